# Update the TaskWorkload column (E) so the decimal values go from a
# single trailing zero ("1,0" / "2,0") to five trailing zeros
# ("1,00000" / "2,00000"), matching the author's "dot to comma prompt
# improvement" formatting change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 5)   # column E = 5
    $text = [string]$cell.Text
    if ([string]::IsNullOrEmpty($text)) {
        continue
    }
    if ($text -match '^(\d+),0$') {
        $cell.Value = $matches[1] + ',00000'
    }
}
